# Applies the coin price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.909.28'
$ws.Range('E2').Value = '  -3.51%  '
$ws.Range('D3').Value = '2.288.58'
$ws.Range('E3').Value = '  -4.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '533.56'
$ws.Range("D5").Style = "Normal"
$ws.Range('E5').Value = '  -4.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '130.58'
$ws.Range("D6").Style = "Normal"
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.581'
$ws.Range("D8").Style = "Normal"
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '2.287.23'
$ws.Range('E9').Value = '  -3.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '0.0994'
$ws.Range("D10").Style = "Normal"
$ws.Range('E10').Value = '  -6.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '5.42'
$ws.Range("D11").Style = "Normal"
$ws.Range('E11').Value = '  -4.42%  '
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('E13').Value = '  -3.96%  '
$ws.Range('E14').Value = '  -4.44%  '
$ws.Range('D15').Value = '2.696.54'
$ws.Range('E15').Value = '  -3.91%  '
$ws.Range('D16').Value = '57.800.99'
$ws.Range('E16').Value = '  -3.62%  '
$ws.Range('E17').Value = '  -4.86%  '
$ws.Range('D18').Value = '2.291.30'
$ws.Range('E18').Value = '  -4.05%  '
$ws.Range('E19').Value = '  -5.58%  '
$ws.Range('E20').Value = '  -5.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '311.95'
$ws.Range("D21").Style = "Normal"
$ws.Range('E21').Value = '  -2.91%  '
$ws.Range('E22').Value = '  -4.81%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '62.35'
$ws.Range("D24").Style = "Normal"
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '0.167'
$ws.Range("D25").Style = "Normal"
$ws.Range('E25').Value = '  -3.89%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '8.01'
$ws.Range("D27").Style = "Normal"
$ws.Range('E27').Value = '  -5.36%  '
$ws.Range('E28').Value = '  -7.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '170.92'
$ws.Range("D29").Style = "Normal"
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('E30').Value = '  -5.90%  '
$ws.Range('D31').Value = '0.0₃0716'
$ws.Range('E31').Value = '  -5.98%  '
$ws.Range('E32').Value = '  -5.40%  '
$ws.Range('E33').Value = '  -6.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '0.378'
$ws.Range("D34").Style = "Normal"
$ws.Range('E34').Value = '  -5.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '17.68'
$ws.Range("D36").Style = "Normal"
$ws.Range('E36').Value = '  -2.64%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  -7.44%  '
$ws.Range('E39').Value = '  -6.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '38.14'
$ws.Range("D40").Style = "Normal"
$ws.Range('E40').Value = '  -1.34%  '
$ws.Range('E41').Value = '  -6.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '141.47'
$ws.Range("D42").Style = "Normal"
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '287.00'
$ws.Range("D43").Style = "Normal"
$ws.Range('E43').Value = '  -10.13%  '
$ws.Range('E44').Value = '  -4.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '0.0946'
$ws.Range("D45").Style = "Normal"
$ws.Range('E45').Value = '  -2.46%  '
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '0.554'
$ws.Range("D47").Style = "Normal"
$ws.Range('E47').Value = '  -3.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '18.04'
$ws.Range("D48").Style = "Normal"
$ws.Range('E48').Value = '  -8.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '0.0209'
$ws.Range("D49").Style = "Normal"
$ws.Range('E49').Value = '  -3.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '10.93'
$ws.Range("D50").Style = "Normal"
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('D51').Value = '0.0₆0200'
$ws.Range('E51').Value = '  +84.61%  '
